$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("GSM")
$ws6 = $wb.Worksheets.Item("LTW")
$ws7 = $wb.Worksheets.Item("WVR")

# ALC
$ws1.Range("H41").Value = 536.6111
$ws1.Range("I41").Value = 272.5
$ws1.Range("J41").Value = 747.9
$ws1.Range("K41").Value = 272.5
$ws1.Range("L41").Value = 747.9
$ws1.Range("M41").Value = 167.5
$ws1.Range("N41").Value = -1627.9
$ws1.Range("H62").Value = 2343.4443
$ws1.Range("I62").Value = 2448.125
$ws1.Range("J62").Value = 1506
$ws1.Range("K62").Value = 2448.125
$ws1.Range("L62").Value = 1506
$ws1.Range("M62").Value = -1824.125
$ws1.Range("N62").Value = -2754
$ws1.Range("H65").Value = 2343.4443
$ws1.Range("I65").Value = 2448.125
$ws1.Range("J65").Value = 1506
$ws1.Range("K65").Value = 12240.625
$ws1.Range("L65").Value = 7530
$ws1.Range("M65").Value = -9120.625
$ws1.Range("N65").Value = -13770
$ws1.Range("H69").Value = 3636.3333
$ws1.Range("I69").Value = 4804.3335
$ws1.Range("K69").Value = 14413.0005
$ws1.Range("M69").Value = -13539.0005
$ws1.Range("H72").Value = 3636.3333
$ws1.Range("I72").Value = 4804.3335
$ws1.Range("K72").Value = 43239.0015
$ws1.Range("M72").Value = -38871.0015
$ws1.Range("H113").Value = 13600
$ws1.Range("I113").Value = 3000
$ws1.Range("J113").Value = 20666.666
$ws1.Range("K113").Value = 3000
$ws1.Range("L113").Value = 20666.666
$ws1.Range("M113").Value = 254
$ws1.Range("N113").Value = -27174.666
$ws1.Range("H137").Value = 4447.952
$ws1.Range("I137").Value = 4718.8335
$ws1.Range("J137").Value = 4086.7778
$ws1.Range("K137").Value = 14156.5005
$ws1.Range("L137").Value = 12260.3334
$ws1.Range("M137").Value = -11606.5005
$ws1.Range("N137").Value = -17360.3334
$ws1.Range("H139").Value = 39175.555
$ws1.Range("J139").Value = 39175.555
$ws1.Range("L139").Value = 39175.555
$ws1.Range("N139").Value = -49455.555

# ARM
$ws2.Range("H2").Value = 1242.3636
$ws2.Range("I2").Value = 1272.8
$ws2.Range("J2").Value = 1177.1428
$ws2.Range("K2").Value = 1272.8
$ws2.Range("L2").Value = 1177.1428
$ws2.Range("M2").Value = -1159.8
$ws2.Range("N2").Value = -1403.1428
$ws2.Range("H97").Value = 1066.2858
$ws2.Range("I97").Value = 910.6667
$ws2.Range("K97").Value = 910.6667
$ws2.Range("M97").Value = -414.6667
$ws2.Range("H116").Value = 1242.3636
$ws2.Range("I116").Value = 1272.8
$ws2.Range("J116").Value = 1177.1428
$ws2.Range("K116").Value = 1272.8
$ws2.Range("L116").Value = 1177.1428
$ws2.Range("M116").Value = 1021.2
$ws2.Range("N116").Value = -5765.1428
$ws2.Range("H122").Value = 2782
$ws2.Range("I122").Value = 1806.8572
$ws2.Range("K122").Value = 5420.571599999999
$ws2.Range("M122").Value = -2970.571599999999

# BSM
$ws3.Range("H3").Value = 1242.3636
$ws3.Range("I3").Value = 1272.8
$ws3.Range("J3").Value = 1177.1428
$ws3.Range("K3").Value = 1272.8
$ws3.Range("L3").Value = 1177.1428
$ws3.Range("M3").Value = -1158.8
$ws3.Range("N3").Value = -1405.1428
$ws3.Range("H122").Value = 0
$ws3.Range("J122").Value = 0
$ws3.Range("L122").Value = 0
$ws3.Range("N122").ClearContents()

# CRP
$ws4.Range("H31").Value = 4727.811
$ws4.Range("I31").Value = 2037.6666
$ws4.Range("J31").Value = 6019.08
$ws4.Range("K31").Value = 2037.6666
$ws4.Range("L31").Value = 6019.08
$ws4.Range("M31").Value = -1742.6666
$ws4.Range("N31").Value = -6609.08
$ws4.Range("H34").Value = 4727.811
$ws4.Range("I34").Value = 2037.6666
$ws4.Range("J34").Value = 6019.08
$ws4.Range("K34").Value = 2037.6666
$ws4.Range("L34").Value = 6019.08
$ws4.Range("M34").Value = -1835.6666
$ws4.Range("N34").Value = -6423.08
$ws4.Range("H107").Value = 540.25714
$ws4.Range("I107").Value = 491.81482
$ws4.Range("J107").Value = 703.75
$ws4.Range("K107").Value = 491.81482
$ws4.Range("L107").Value = 703.75
$ws4.Range("M107").Value = 1428.18518
$ws4.Range("N107").Value = -4543.75
$ws4.Range("H122").Value = 4683.3335
$ws4.Range("I122").Value = 4950
$ws4.Range("J122").Value = 4630
$ws4.Range("K122").Value = 14850
$ws4.Range("L122").Value = 13890
$ws4.Range("M122").Value = -12400
$ws4.Range("N122").Value = -18790

# GSM
$ws5.Range("H93").Value = 17499.875
$ws5.Range("J93").Value = 17499.875
$ws5.Range("L93").Value = 17499.875
$ws5.Range("N93").Value = -21243.875

# LTW
$ws6.Range("H22").Value = 2542.524
$ws6.Range("I22").Value = 2090.0908
$ws6.Range("J22").Value = 3040.2
$ws6.Range("K22").Value = 2090.0908
$ws6.Range("L22").Value = 3040.2
$ws6.Range("M22").Value = -1795.0908
$ws6.Range("N22").Value = -3630.2
$ws6.Range("H27").Value = 2542.524
$ws6.Range("I27").Value = 2090.0908
$ws6.Range("J27").Value = 3040.2
$ws6.Range("K27").Value = 2090.0908
$ws6.Range("L27").Value = 3040.2
$ws6.Range("M27").Value = -1983.0908
$ws6.Range("N27").Value = -3254.2
$ws6.Range("H122").Value = 4259.222
$ws6.Range("I122").Value = 2698.8333
$ws6.Range("K122").Value = 8096.499899999999
$ws6.Range("M122").Value = -5646.499899999999

# WVR
$ws7.Range("H62").Value = 38559310
$ws7.Range("I62").Value = 125003740
$ws7.Range("J62").Value = 139566
$ws7.Range("K62").Value = 125003740
$ws7.Range("L62").Value = 139566
$ws7.Range("M62").Value = -125003116
$ws7.Range("N62").Value = -140814
$ws7.Range("H65").Value = 38559310
$ws7.Range("I65").Value = 125003740
$ws7.Range("J65").Value = 139566
$ws7.Range("K65").Value = 625018700
$ws7.Range("L65").Value = 697830
$ws7.Range("M65").Value = -625015580
$ws7.Range("N65").Value = -704070
$ws7.Range("H81").Value = 17858124
$ws7.Range("I81").Value = 20090292
$ws7.Range("J81").Value = 787.5
$ws7.Range("K81").Value = 40180584
$ws7.Range("L81").Value = 1575
$ws7.Range("M81").Value = -40179523
$ws7.Range("N81").Value = -3697
$ws7.Range("H84").Value = 17858124
$ws7.Range("I84").Value = 20090292
$ws7.Range("J84").Value = 787.5
$ws7.Range("K84").Value = 200902920
$ws7.Range("L84").Value = 7875
$ws7.Range("M84").Value = -200897616
$ws7.Range("N84").Value = -18483
$ws7.Range("H86").Value = 28500
$ws7.Range("J86").Value = 28500
$ws7.Range("L86").Value = 28500
$ws7.Range("N86").Value = -30746
$ws7.Range("H89").Value = 28500
$ws7.Range("J89").Value = 28500
$ws7.Range("L89").Value = 142500
$ws7.Range("N89").Value = -153732
$ws7.Range("H96").Value = 128139370
$ws7.Range("I96").Value = 500250000
$ws7.Range("J96").Value = 4102487.5
$ws7.Range("K96").Value = 500250000
$ws7.Range("L96").Value = 4102487.5
$ws7.Range("M96").Value = -500248627
$ws7.Range("N96").Value = -4105233.5
$ws7.Range("H122").Value = 7696.4165
$ws7.Range("I122").Value = 6489.4
$ws7.Range("K122").Value = 19468.2
$ws7.Range("M122").Value = -17018.2
